$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content change ---------------------------------------------------
# "Journée Sport et Détente à Beaucourt (à confirmer)" used to sit at row 22
# (Dim 27 Juin). It has been removed from there (rows 23-29 shift up one),
# and re-added as a new row 29 with an updated date (Dim 4 Juillet) and a
# trimmed title (the "(à confirmer)" suffix is gone, now it's confirmed).
$ws.Rows(22).Delete()

$ws.Rows(29).Insert()
$ws.Range("A29").Value = "Dim 4 Juillet"
$ws.Range("B29").Value = "Journée Sport et Détente à Beaucourt"
$ws.Range("C29").Value = "Beaucourt OS"
$ws.Range("D29").Value = "Randonnée"
$ws.Range("E29").Value = "rando_beaucourt"

# --- View / selection change -------------------------------------------
# The sheet's visible selection moved from B58 to B30 (scrolled up, following
# the shifted rows).
$ws.Range("B30").Select()
